$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 11, shifting existing rows 11-18 down to 12-19
$ws.Rows("11:11").Insert()

# The freshly inserted row doesn't inherit the bordered data-row formatting;
# copy it from the row below (the original row 11, now shifted to row 12).
$ws.Range("A12:H12").Copy() | Out-Null
$ws.Range("A11:H11").PasteSpecial(-4122) | Out-Null

# Capture existing hyperlinks (their positions did not auto-shift) so we can
# rebuild them pointing at the correct (shifted) cells, preserving order.
$links = @()
foreach ($hl in $ws.Hyperlinks) {
    $links += ,@($hl.Range.Row, $hl.Range.Column, $hl.Address)
}
$ws.Hyperlinks.Delete()
foreach ($l in $links) {
    $r = $l[0]
    if ($r -ge 11) { $r = $r + 1 }
    $c = $l[1]
    $addr = $l[2]
    $ws.Hyperlinks.Add($ws.Cells.Item($r, $c), $addr) | Out-Null
}

# Populate the new row 11 with the 0 ohm jumper part
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 2
$ws.Range("C11").Value = "0.0QBK-ND"
$ws.Range("D11").Value = "ZOR-25-B-52-0R"
$ws.Range("E11").Value = "Yageo"
$ws.Range("F11").Value = "0 OHM 1/4W JUMP"
$ws.Range("G11").Value = 0.1
$ws.Range("H11").Formula = "=G11*B11"

# Update index numbers for shifted rows (now 12, 13, 14)
$ws.Range("A12").Value = 11
$ws.Range("A13").Value = 12
$ws.Range("A14").Value = 13

# Add hyperlink for the new row's Digi-Key part number (added last -> new rId)
$ws.Hyperlinks.Add($ws.Range("C11"), "https://www.digikey.com/products/en?keywords=ZOR-25-B-52-0R") | Out-Null
